# revisi menrevisi kombinasi loc
# Recomputed face-recognition match results (distance/threshold/confidence/
# recall columns B:E, predicted-name F, correctness G) for every probe image
# in A1:G45 -- the row set also gained one more "TD" probe (TD_4.png), so the
# previous TO_1..TD_4 block shifts down by one row versus the old sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=image, B=distance, C=hamming/threshold, D=confidence, E=recall, F=predicted name, G=Benar/Salah
$rows = @(
    ,@("AK_1.png", 0.55, 0.002, 0.971, 0.267, "Tidak Diketahui", "Salah")
    ,@("AK_2.png", 0.904, 0.003, 0.982, 0.3, "Tidak Diketahui", "Salah")
    ,@("AK_3.png", 0.89, 0.003, 0.973, 0.533, "Akhlak Kamiswara", "Benar")
    ,@("AK_4.png", 0.548, 0.002, 0.962, 0.4, "Tidak Diketahui", "Salah")
    ,@("AK_5.png", 0.536, 0.002, 0.984, 0.267, "Tidak Diketahui", "Salah")
    ,@("MIB_1.png", 1.619, 0.005, 0.971, 0.533, "Muhammad Iqbal Baqi", "Benar")
    ,@("MIB_2.png", 1.155, 0.004, 0.965, 0.467, "Tidak Diketahui", "Salah")
    ,@("MIB_3.png", 1.8, 0.006, 0.996, 0.9330000000000001, "Muhammad Iqbal Baqi", "Benar")
    ,@("MIB_4.png", 1.165, 0.004, 0.979, 0.6, "Muhammad Iqbal Baqi", "Benar")
    ,@("MIB_5.png", 1.569, 0.005, 0.977, 0.5669999999999999, "Muhammad Iqbal Baqi", "Benar")
    ,@("AAH_1.png", 0.876, 0.003, 0.971, 0.433, "Tidak Diketahui", "Salah")
    ,@("AAH_2.png", 1.109, 0.004, 0.988, 0.7, "Andrea Ayunove Hutami", "Benar")
    ,@("AAH_3.png", 0.906, 0.003, 0.967, 0.5669999999999999, "Andrea Ayunove Hutami", "Benar")
    ,@("TI_1.png", 0.956, 0.003, 0.9419999999999999, 0.367, "Tidak Diketahui", "Salah")
    ,@("TI_2.png", 0.95, 0.003, 0.983, 0.3, "Tidak Diketahui", "Salah")
    ,@("TI_3.png", 0.758, 0.003, 0.985, 0.333, "Tidak Diketahui", "Salah")
    ,@("TI_4.png", 0.744, 0.002, 0.958, 0.467, "Tidak Diketahui", "Salah")
    ,@("TI_5.png", 1.047, 0.003, 0.985, 0.467, "Tidak Diketahui", "Salah")
    ,@("RAS_1.png", 0.6870000000000001, 0.002, 0.97, 0.2, "Tidak Diketahui", "Salah")
    ,@("RAS_2.png", 1.096, 0.004, 0.96, 0.4, "Tidak Diketahui", "Salah")
    ,@("RAS_3.png", 0.578, 0.002, 0.978, 0.233, "Tidak Diketahui", "Salah")
    ,@("RAS_4.png", 1.29, 0.004, 0.967, 0.233, "Tidak Diketahui", "Salah")
    ,@("RAS_5.png", 1.2, 0.004, 0.965, 0.5, "Ridha Ayu Salsabila", "Benar")
    ,@("RR_1.png", 1.488, 0.005, 0.976, 0.3, "Tidak Diketahui", "Salah")
    ,@("RR_2.png", 1.364, 0.005, 0.981, 0.333, "Tidak Diketahui", "Salah")
    ,@("RR_3.png", 1.141, 0.004, 0.895, 0.3, "Tidak Diketahui", "Salah")
    ,@("RR_4.png", 1.407, 0.005, 0.986, 0.333, "Tidak Diketahui", "Salah")
    ,@("RR_5.png", 1.484, 0.005, 0.979, 0.5, "Rafiqo Rapitasari", "Benar")
    ,@("AR_1.png", 0.889, 0.003, 0.955, 0.5, "Arizli Romadhon", "Benar")
    ,@("GA_1.png", 1.706, 0.006, 0.981, 0.633, "Gege Ardiyansyah", "Benar")
    ,@("GA_2.png", 0.68, 0.002, 0.981, 0.267, "Tidak Diketahui", "Salah")
    ,@("GA_3.png", 0.777, 0.003, 0.981, 0.3, "Tidak Diketahui", "Salah")
    ,@("FY_1.png", 1.165, 0.004, 0.964, 0.233, "Tidak Diketahui", "Salah")
    ,@("FY_2.png", 1.695, 0.006, 0.969, 0.2, "Tidak Diketahui", "Salah")
    ,@("FY_3.png", 1.458, 0.005, 0.984, 0.267, "Tidak Diketahui", "Salah")
    ,@("FY_4.png", 1.268, 0.004, 0.964, 0.367, "Tidak Diketahui", "Salah")
    ,@("TO_1.png", 0.803, 0.003, 0.962, 0.333, "Tidak Diketahui", "Salah")
    ,@("TO_2.png", 1.113, 0.004, 0.984, 0.333, "Tidak Diketahui", "Salah")
    ,@("TO_3.png", 0.872, 0.003, 0.977, 0.4, "Tidak Diketahui", "Salah")
    ,@("TO_4.png", 7.452, 0.025, 0.925, 0.367, "Tidak Diketahui", "Salah")
    ,@("TO_5.png", 7.112, 0.023, 0.925, 0.267, "Tidak Diketahui", "Salah")
    ,@("TD_1.png", 3.8, 0.013, 0.885, 0.433, "Tidak Diketahui", "Benar")
    ,@("TD_2.png", 5.056, 0.017, 0.873, 0.267, "Tidak Diketahui", "Benar")
    ,@("TD_3.png", 1.061, 0.004, 0.9409999999999999, 0.267, "Tidak Diketahui", "Benar")
    ,@("TD_4.png", 0.995, 0.003, 0.979, 0.233, "Tidak Diketahui", "Benar")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
